$d = $word.ActiveDocument

# Bug #449: two stray empty paragraphs were left in the generated
# document right after the "While ... are reached:" block and right
# before the numbered list ("list item 1," etc.). Both paragraphs
# contain no text at all -- only a leftover "_GoBack" bookmark -- and
# are solely distinguished by their left indentation (560 and 1120
# twips respectively). Remove both of them.
$found = 0
$i = $d.Paragraphs.Count
while ($i -ge 1) {
    $p = $d.Paragraphs($i)
    $indent = $p.Format.LeftIndent
    $text = $p.Range.Text
    $isEmpty = ($text -eq "" -or $text -eq [char]13 -or $text -eq [char]7)
    if ($isEmpty -and ($indent -eq 28 -or $indent -eq 56)) {
        $p.Range.Delete()
        $found = $found + 1
    }
    $i = $i - 1
}
